$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D/E columns stay text (avoid Excel auto-converting numeric-looking strings)
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.275.06'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.26%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.344.14'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.89%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.18'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +5.32%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '185.49'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.21%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.85%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.14%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '46.97'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.99%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.50%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '669.83'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +10.88%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.876.50'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.88%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.49'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.78%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.442.81'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.46%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.88'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.69%  '
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.341.25'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.88%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.10'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.22%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.897'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.59%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -4.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '101.30'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.04'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.84%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.03'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.28%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.78'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.29%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.96%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '32.28'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +6.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.51'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.63%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.78'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.30%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '614.63'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +6.01%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +2.58%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.08'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.60%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.866.45'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.84%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.25%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.15%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '56.28'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.75%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.30%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0₃0702'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -3.12%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.56%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '32.86'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.46%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.72%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.42'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.98%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.65%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.35%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -16.13%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.40%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.23%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.54'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.95%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.31'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +3.75%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '129.12'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +4.83%  '
